$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as plain text, avoiding Excel
# auto-conversion of numeric-looking strings into floating point numbers,
# and then reset the cell style so no extra "Text" number format sticks.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '62.634.56'
$ws.Range("E2").Value = '  -4.27%  '

Set-TextValue $ws.Range("D3") '3.273.67'
$ws.Range("E3").Value = '  -6.19%  '

$ws.Range("E4").Value = '  +0.04%  '

Set-TextValue $ws.Range("D5") '538.06'
$ws.Range("E5").Value = '  -2.82%  '

Set-TextValue $ws.Range("D6") '169.91'
$ws.Range("E6").Value = '  -4.70%  '

$ws.Range("E7").Value = '  -4.67%  '

Set-TextValue $ws.Range("D9") '3.266.07'
$ws.Range("E9").Value = '  -6.17%  '

Set-TextValue $ws.Range("D10") '0.604'
$ws.Range("E10").Value = '  -4.31%  '

Set-TextValue $ws.Range("D11") '0.150'
$ws.Range("E11").Value = '  -1.40%  '

Set-TextValue $ws.Range("D12") '52.29'
$ws.Range("E12").Value = '  -2.29%  '

$ws.Range("E13").Value = '  -3.24%  '

Set-TextValue $ws.Range("D14") '8.75'
$ws.Range("E14").Value = '  -5.03%  '

Set-TextValue $ws.Range("D15") '3.791.54'
$ws.Range("E15").Value = '  -6.23%  '

Set-TextValue $ws.Range("D16") '17.73'
$ws.Range("E16").Value = '  -3.82%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D17") '0.116'
$ws.Range("E17").Value = '  -4.43%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D18") '3.278.19'
$ws.Range("E18").Value = '  -5.97%  '

$ws.Range("E19").Value = '  -4.63%  '

Set-TextValue $ws.Range("D20") '62.566.47'
$ws.Range("E20").Value = '  -4.35%  '

Set-TextValue $ws.Range("D21") '0.960'
$ws.Range("E21").Value = '  -2.97%  '

Set-TextValue $ws.Range("D22") '410.41'
$ws.Range("E22").Value = '  -1.58%  '

Set-TextValue $ws.Range("D23") '4.33'
$ws.Range("E23").Value = '  +5.25%  '

$ws.Range("E24").Value = '  -1.49%  '

Set-TextValue $ws.Range("D25") '13.42'
$ws.Range("E25").Value = '  +5.63%  '

Set-TextValue $ws.Range("D26") '82.06'
$ws.Range("E26").Value = '  -4.90%  '

Set-TextValue $ws.Range("D27") '10.39'
$ws.Range("E27").Value = '  -3.42%  '

$ws.Range("E28").Value = '  -5.41%  '

Set-TextValue $ws.Range("D29") '8.48'
$ws.Range("E29").Value = '  -5.58%  '

Set-TextValue $ws.Range("D30") '28.66'
$ws.Range("E30").Value = '  -5.17%  '

Set-TextValue $ws.Range("D31") '6.26'
$ws.Range("E31").Value = '  -3.40%  '

Set-TextValue $ws.Range("D32") '11.20'
$ws.Range("E32").Value = '  -4.19%  '

Set-TextValue $ws.Range("D33") '569.65'
$ws.Range("E33").Value = '  -6.43%  '

$ws.Range("E34").Value = '  -4.23%  '

Set-TextValue $ws.Range("D35") '57.47'
$ws.Range("E35").Value = '  -3.37%  '

$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("E37").Value = '  -1.72%  '

Set-TextValue $ws.Range("D38") '34.62'
$ws.Range("E38").Value = '  -7.12%  '

Set-TextValue $ws.Range("D39") '3.37'
$ws.Range("E39").Value = '  +3.90%  '

Set-TextValue $ws.Range("D40") '0.0₃0725'
$ws.Range("E40").Value = '  -7.67%  '

Set-TextValue $ws.Range("D41") '0.361'
$ws.Range("E41").Value = '  -5.00%  '

Set-TextValue $ws.Range("D42") '3.089.48'
$ws.Range("E42").Value = '  -8.49%  '

Set-TextValue $ws.Range("D43") '0.997'
$ws.Range("E43").Value = '  -0.16%  '

Set-TextValue $ws.Range("D44") '3.20'
$ws.Range("E44").Value = '  -1.26%  '

$ws.Range("E45").Value = '  -4.20%  '

Set-TextValue $ws.Range("D46") '0.0395'
$ws.Range("E46").Value = '  -4.31%  '

$ws.Range("E47").Value = '  -5.94%  '

$ws.Range("E48").Value = '  -4.63%  '

$ws.Range("E49").Value = '  -3.98%  '

Set-TextValue $ws.Range("D50") '131.68'
$ws.Range("E50").Value = '  -4.42%  '

Set-TextValue $ws.Range("D51") '7.91'
$ws.Range("E51").Value = '  -6.66%  '
